$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

$t.Cell(1, 1).Range.Text = "27 x 81" + $nl + "  8    1" + $nl + "  ----" + $nl + "2|    |" + $nl + "7|    |"
$t.Cell(1, 2).Range.Text = "71 x 84" + $nl + "  8    4" + $nl + "  ----" + $nl + "7|    |" + $nl + "1|    |"
$t.Cell(1, 3).Range.Text = "70 x 43" + $nl + "  4    3" + $nl + "  ----" + $nl + "7|    |" + $nl + "0|    |"
$t.Cell(2, 1).Range.Text = "53 x 35" + $nl + "  3    5" + $nl + "  ----" + $nl + "5|    |" + $nl + "3|    |"
$t.Cell(2, 2).Range.Text = "32 x 19" + $nl + "  1    9" + $nl + "  ----" + $nl + "3|    |" + $nl + "2|    |"
$t.Cell(2, 3).Range.Text = "51 x 41" + $nl + "  4    1" + $nl + "  ----" + $nl + "5|    |" + $nl + "1|    |"
$t.Cell(3, 1).Range.Text = "99 x 78" + $nl + "  7    8" + $nl + "  ----" + $nl + "9|    |" + $nl + "9|    |"
$t.Cell(3, 2).Range.Text = "28 x 24" + $nl + "  2    4" + $nl + "  ----" + $nl + "2|    |" + $nl + "8|    |"
$t.Cell(3, 3).Range.Text = "24 x 49" + $nl + "  4    9" + $nl + "  ----" + $nl + "2|    |" + $nl + "4|    |"
$t.Cell(4, 1).Range.Text = "29 x 38" + $nl + "  3    8" + $nl + "  ----" + $nl + "2|    |" + $nl + "9|    |"
$t.Cell(4, 2).Range.Text = "70 x 17" + $nl + "  1    7" + $nl + "  ----" + $nl + "7|    |" + $nl + "0|    |"
$t.Cell(4, 3).Range.Text = "24 x 46" + $nl + "  4    6" + $nl + "  ----" + $nl + "2|    |" + $nl + "4|    |"
$t.Cell(5, 1).Range.Text = "50 x 90" + $nl + "  9    0" + $nl + "  ----" + $nl + "5|    |" + $nl + "0|    |"
$t.Cell(5, 2).Range.Text = "81 x 77" + $nl + "  7    7" + $nl + "  ----" + $nl + "8|    |" + $nl + "1|    |"
$t.Cell(5, 3).Range.Text = "59 x 93" + $nl + "  9    3" + $nl + "  ----" + $nl + "5|    |" + $nl + "9|    |"
